$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift each data row one column to the left (dropping the oldest/leftmost
# error estimate). Rows 2-10 gain a brand new value in column K (the newly
# evaluated ifoCAST data point), while rows 11-20 simply lose their
# right-most (now out-of-range) value as the staircase shrinks.

# Row 2
$ws.Range("B2").Value = -1.876881696346878
$ws.Range("C2").Value = -1.938055874761079
$ws.Range("D2").Value = 0.1034000623951754
$ws.Range("E2").Value = -0.4791944746384656
$ws.Range("F2").Value = -0.3835194338291078
$ws.Range("G2").Value = -0.6439384708278306
$ws.Range("H2").Value = 0.1915165253306606
$ws.Range("I2").Value = 0.09143716542803909
$ws.Range("J2").Value = -0.2623493406516572
$ws.Range("K2").Value = -0.1141719502192118

# Row 3
$ws.Range("B3").Value = -2.04667040697722
$ws.Range("C3").Value = -0.005214469820965406
$ws.Range("D3").Value = -0.5878090068546065
$ws.Range("E3").Value = -0.4921339660452486
$ws.Range("F3").Value = -0.7525530030439714
$ws.Range("G3").Value = 0.08290199311451979
$ws.Range("H3").Value = -0.01717736678810172
$ws.Range("I3").Value = -0.370963872867798
$ws.Range("J3").Value = -0.2227864824353526
$ws.Range("K3").Value = -0.06534049343119763

# Row 4
$ws.Range("B4").Value = 0.2028872511467736
$ws.Range("C4").Value = -0.3797072858868674
$ws.Range("D4").Value = -0.2840322450775096
$ws.Range("E4").Value = -0.5444512820762324
$ws.Range("F4").Value = 0.2910037140822588
$ws.Range("G4").Value = 0.1909243541796373
$ws.Range("H4").Value = -0.1628621519000589
$ws.Range("I4").Value = -0.01468476146761361
$ws.Range("J4").Value = 0.1427612275365414
$ws.Range("K4").Value = -0.3578982289895294

# Row 5
$ws.Range("B5").Value = -0.263209115513067
$ws.Range("C5").Value = -0.1675340747037092
$ws.Range("D5").Value = -0.427953111702432
$ws.Range("E5").Value = 0.4075018844560592
$ws.Range("F5").Value = 0.3074225245534377
$ws.Range("G5").Value = -0.04636398152625851
$ws.Range("H5").Value = 0.1018134089061868
$ws.Range("I5").Value = 0.2592593979103418
$ws.Range("J5").Value = -0.241400058615729
$ws.Range("K5").Value = 0.03895782191974018

# Row 6
$ws.Range("B6").Value = -0.1220406607065741
$ws.Range("C6").Value = -0.3824596977052969
$ws.Range("D6").Value = 0.4529952984531944
$ws.Range("E6").Value = 0.3529159385505728
$ws.Range("F6").Value = -0.0008705675291234075
$ws.Range("G6").Value = 0.1473068229033219
$ws.Range("H6").Value = 0.3047528119074769
$ws.Range("I6").Value = -0.1959066446185939
$ws.Range("J6").Value = 0.08445123591687528
$ws.Range("K6").Value = -0.08752362192813637

# Row 7
$ws.Range("B7").Value = -0.4000138190652613
$ws.Range("C7").Value = 0.4354411770932299
$ws.Range("D7").Value = 0.3353618171906084
$ws.Range("E7").Value = -0.01842468888908786
$ws.Range("F7").Value = 0.1297527015433575
$ws.Range("G7").Value = 0.2871986905475125
$ws.Range("H7").Value = -0.2134607659785583
$ws.Range("I7").Value = 0.06689711455691082
$ws.Range("J7").Value = -0.1050777432881008
$ws.Range("K7").Value = 0.2007847964582307

# Row 8
$ws.Range("B8").Value = 0.5235484960442361
$ws.Range("C8").Value = 0.4234691361416146
$ws.Range("D8").Value = 0.06968263006191837
$ws.Range("E8").Value = 0.2178600204943637
$ws.Range("F8").Value = 0.3753060094985187
$ws.Range("G8").Value = -0.1253534470275521
$ws.Range("H8").Value = 0.155004433507917
$ws.Range("I8").Value = -0.01697042433709459
$ws.Range("J8").Value = 0.2888921154092369
$ws.Range("K8").Value = -0.3234378372780493

# Row 9
$ws.Range("B9").Value = 0.6020437696293982
$ws.Range("C9").Value = 0.248257263549702
$ws.Range("D9").Value = 0.3964346539821473
$ws.Range("E9").Value = 0.5538806429863024
$ws.Range("F9").Value = 0.0532211864602315
$ws.Range("G9").Value = 0.3335790669957007
$ws.Range("H9").Value = 0.161604209150689
$ws.Range("I9").Value = 0.4674667488970205
$ws.Range("J9").Value = -0.1448632037902657
$ws.Range("K9").Value = 0.5435081814088458

# Row 10
$ws.Range("B10").Value = 1.160284491302352
$ws.Range("C10").Value = 1.308461881734797
$ws.Range("D10").Value = 1.465907870738952
$ws.Range("E10").Value = 0.9652484142128814
$ws.Range("F10").Value = 1.245606294748351
$ws.Range("G10").Value = 1.073631436903339
$ws.Range("H10").Value = 1.379493976649671
$ws.Range("I10").Value = 0.7671640239623843
$ws.Range("J10").Value = 1.455535409161496
$ws.Range("K10").Value = 1.178844253737389

# Row 11
$ws.Range("B11").Value = 0.3644876457690405
$ws.Range("C11").Value = 0.5219336347731955
$ws.Range("D11").Value = 0.02127417824712469
$ws.Range("E11").Value = 0.3016320587825939
$ws.Range("F11").Value = 0.1296572009375822
$ws.Range("G11").Value = 0.4355197406839137
$ws.Range("H11").Value = -0.1768102120033725
$ws.Range("I11").Value = 0.511561173195739
$ws.Range("J11").Value = 0.2348700177716323
$ws.Range("K11").ClearContents()

# Row 12
$ws.Range("B12").Value = 0.5259015322863045
$ws.Range("C12").Value = 0.0252420757602338
$ws.Range("D12").Value = 0.305599956295703
$ws.Range("E12").Value = 0.1336250984506913
$ws.Range("F12").Value = 0.4394876381970228
$ws.Range("G12").Value = -0.1728423144902634
$ws.Range("H12").Value = 0.5155290707088481
$ws.Range("I12").Value = 0.2388379152847414
$ws.Range("J12").ClearContents()

# Row 13
$ws.Range("B13").Value = 0.1608821659304752
$ws.Range("C13").Value = 0.4412400464659443
$ws.Range("D13").Value = 0.2692651886209327
$ws.Range("E13").Value = 0.5751277283672642
$ws.Range("F13").Value = -0.03720222432002201
$ws.Range("G13").Value = 0.6511691608790895
$ws.Range("H13").Value = 0.3744780054549828
$ws.Range("I13").ClearContents()

# Row 14
$ws.Range("B14").Value = 0.2004338646102796
$ws.Range("C14").Value = 0.028459006765268
$ws.Range("D14").Value = 0.3343215465115995
$ws.Range("E14").Value = -0.2780084061756867
$ws.Range("F14").Value = 0.4103629790234248
$ws.Range("G14").Value = 0.1336718235993181
$ws.Range("H14").ClearContents()

# Row 15
$ws.Range("B15").Value = -0.01687220848682837
$ws.Range("C15").Value = 0.2889903312595031
$ws.Range("D15").Value = -0.3233396214277831
$ws.Range("E15").Value = 0.3650317637713285
$ws.Range("F15").Value = 0.08834060834722172
$ws.Range("G15").ClearContents()

# Row 16
$ws.Range("B16").Value = 0.2221289093234493
$ws.Range("C16").Value = -0.3902010433638369
$ws.Range("D16").Value = 0.2981703418352746
$ws.Range("E16").Value = 0.02147918641116785
$ws.Range("F16").ClearContents()

# Row 17
$ws.Range("B17").Value = -0.4197872457205535
$ws.Range("C17").Value = 0.268584139478558
$ws.Range("D17").Value = -0.00810701594554874
$ws.Range("E17").ClearContents()

# Row 18
$ws.Range("B18").Value = 0.2504334827489171
$ws.Range("C18").Value = -0.02625767267518964
$ws.Range("D18").ClearContents()

# Row 19
$ws.Range("B19").Value = -0.04428949692388896
$ws.Range("C19").ClearContents()

$ws.Range("B20").ClearContents()

Write-Host "Applied ifoCAST full series evaluation shift"
